# issue #5: stock data from json to db
#
# The 股票 (stock) sheet is extended with the "category", "source_file"
# and "index" columns that the other json->db-migrated sheets already
# carry:
#   - a new "category" column is inserted right after "property_category"
#     (shifting "date", "legislator_name" and "legislator_id" one column
#     to the right);
#   - "source_file" and "index" columns are appended after
#     "legislator_id".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Capture old values before anything gets overwritten. The date cell
# ("2012-04-27") is plain text that merely looks like a date, so its text
# is re-prefixed with an apostrophe on write-back and then restyled -
# otherwise the COM layer (like real Excel) auto-detects the date-shaped
# string and silently turns it into a date serial number.
$oldDateHdr    = $ws.Cells.Item(1, 9).Value()
$oldLegNameHdr = $ws.Cells.Item(1, 10).Value()
$oldLegIdHdr   = $ws.Cells.Item(1, 11).Value()

$oldDateVal    = $ws.Cells.Item(2, 9).Text
$oldDateStyle  = $ws.Cells.Item(2, 9).Style
$oldLegNameVal = $ws.Cells.Item(2, 10).Value()
$oldLegIdVal   = $ws.Cells.Item(2, 11).Value()

# --- Row 1 (headers) ----------------------------------------------------
$ws.Cells.Item(1, 9).Value  = "category"        # I1 = category (new)
$ws.Cells.Item(1, 10).Value = $oldDateHdr         # J1 = date (was I1)
$ws.Cells.Item(1, 11).Value = $oldLegNameHdr      # K1 = legislator_name (was J1)
$ws.Cells.Item(1, 12).Value = $oldLegIdHdr        # L1 = legislator_id (was K1)
$ws.Cells.Item(1, 13).Value = "source_file"       # M1 = source_file (new)
$ws.Cells.Item(1, 14).Value = "index"             # N1 = index (new)

# Re-apply the header look (bold, centered, thin border) to every header
# cell on the row so the newly created ones match the existing ones.
$hdrRow = $ws.Range($ws.Cells.Item(1, 2), $ws.Cells.Item(1, 14))
$hdrRow.Font.Bold = $true
$hdrRow.HorizontalAlignment = -4108
$hdrRow.VerticalAlignment = -4160
$hdrRow.Borders.LineStyle = 1

# --- Row 2 (data) --------------------------------------------------------
$ws.Cells.Item(2, 9).Value  = "normal"             # I2 = category value
$ws.Cells.Item(2, 10).Value = "'" + $oldDateVal      # J2 = date (was I2)
$ws.Cells.Item(2, 11).Value = $oldLegNameVal         # K2 = legislator_name (was J2)
$ws.Cells.Item(2, 12).Value = $oldLegIdVal           # L2 = legislator_id value (was K2 = 1750)
$ws.Cells.Item(2, 13).Value = "tmp1dd71"             # M2 = source_file value
$ws.Cells.Item(2, 14).Value = 75                     # N2 = index value

# The original "date" cell carried plain (unstyled) formatting; restore it
# after the apostrophe-prefixed write-back above.
$ws.Cells.Item(2, 10).Style = $oldDateStyle
